$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:email3@example.com", "", "", "email3@example.com")
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
Write-Host "added"
